$wb = $excel.ActiveWorkbook

# ---------- Sheet: Detailed Expenses ----------
$ws1 = $wb.Worksheets.Item("Detailed Expenses")

# Row 2: Petrol expense - update date, fix category, add note
$ws1.Range("A2").Value = "29/8/2025"
$ws1.Range("C2").Value = "Petrol"
$ws1.Range("G2").Value = "for car"

# Row 3: Salary income - update date/description, increase amount, add note
$ws1.Range("A3").Value = "29/8/2025"
$ws1.Range("B3").Value = "salary"
$ws1.Range("D3").Value = 1000000
$ws1.Range("G3").Value = "Yearly Salary"

# Row 4: remove old entry entirely
$ws1.Range("A4:G4").Clear()

# Row 5: remove old entry, then write the new TOTAL row
$ws1.Range("A5:G5").Clear()
$ws1.Range("B5").Value = "TOTAL"
$ws1.Range("D5").Value = 1000000
$ws1.Range("E5").Value = 1000
$ws1.Range("F5").Value = 999000

# Row 6: remove old entry entirely
$ws1.Range("A6:G6").Clear()

# Row 8: remove old TOTAL row entirely (total moved to row 5)
$ws1.Range("A8:G8").Clear()

# ---------- Sheet: Yearly Summary ----------
$ws2 = $wb.Worksheets.Item("Yearly Summary")
$ws2.Range("B2").Value = 1000000
$ws2.Range("C2").Value = 1000
$ws2.Range("D2").Value = 999000
$ws2.Range("B4").Value = 1000000
$ws2.Range("C4").Value = 1000
$ws2.Range("D4").Value = 999000

# ---------- Sheet: Yearly Notes ----------
$ws3 = $wb.Worksheets.Item("Yearly Notes")
$ws3.Range("A2").Value = "29/8/2025"
$ws3.Range("B2").Value = "for car"
$ws3.Range("A3").Value = "29/8/2025"
$ws3.Range("B3").Value = "Yearly Salary"
